$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ENCAISSEMENT")

# Row 4 (Espèces) - update B4, C4, D4
$ws.Range("B4").Value = 237464.9
$ws.Range("C4").Value = 6721787.58
$ws.Range("D4").Value = 6521787.58

# Row 9 (TOTAL) - update B9, C9, D9
$ws.Range("B9").Value = 634164.65
$ws.Range("C9").Value = 10589012.89
$ws.Range("D9").Value = 10389012.89
